# Regenerate orders with updated distance/size codes.
# The experiment's distance and size condition labels changed:
#   D80 -> D86, D64 -> D69, D51 -> D55, S30 -> S31
# These codes appear as substrings inside several text columns
# (Condition, Filename_Left, Filename_Right, Distance, Size), so a
# plain text substitution across the used range reproduces the diff
# (which only touches the shared-string table, not the cell layout).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange

# Order matters only in that no replacement's target string may be
# mistaken for another replacement's source string later on; verified
# there is no overlap among {D80,D86,D64,D69,D51,D55,S30,S31}.
$usedRange.Replace("D80", "D86")
$usedRange.Replace("D64", "D69")
$usedRange.Replace("D51", "D55")
$usedRange.Replace("S30", "S31")
